$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1 to make room for a header row. This shifts the
# existing poll-figures row (old row 1) down to row 2, and the percentage
# row (old row 3) down to row 4 - its formulas follow the shift and now
# reference row 2 instead of row 1.
$ws.Rows.Item(1).Insert()

# New header row (row 1): party abbreviations.
$ws.Range("A1").Value = "LNP"
$ws.Range("B1").Value = "ALP"
$ws.Range("C1").Value = "GRN"
$ws.Range("D1").Value = "ONP"
$ws.Range("E1").Value = "NXT"
$ws.Range("F1").Value = "UAP"
$ws.Range("G1").Value = "OTH"

# Updated poll figures, now on row 2 (I2's SUMIF and row 4's percentages
# recalculate automatically since they reference row 2 by formula).
$ws.Range("A2").Value = 37
$ws.Range("B2").Value = 35
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 3
$ws.Range("G2").Value = 7

# Selection now sits on the recalculated percentage row.
$ws.Range("A4:G4").Select()

$wb.Save()
